$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs","Bmp3","Acvr2b","ECs",2,0.6666666666666666,1.720049666666667,5.160149000000001,0.6449101994727137,0.6449101994727137,3,1,0.5292956666666667,1.587887,0.135651968140022,0.1356519681400219,0.9104148350181114,8.193733515163002,0.08748333783204776,0.08748333783204773),
    @("FAPs","Bmp3","Acvr2b","FAPs",2,0.6666666666666666,1.720049666666667,5.160149000000001,0.6449101994727137,0.6449101994727137,3,1,1.362890666666667,4.088672000000001,0.3492921120199358,0.3492921120199358,2.344239636903112,21.09815673212801,0.2252620456370223,0.2252620456370223),
    @("FAPs","Bmp3","Acvr2b","sCs",2,0.6666666666666666,1.720049666666667,5.160149000000001,0.6449101994727137,0.6449101994727137,3,1,2.009678666666666,6.029036,0.5150559198400423,0.5150559198400423,3.456747120707111,31.110724086364,0.3321648160036437,0.3321648160036437),
    @("sCs","Bmp3","Acvr2b","ECs",3,1,0.9470653333333333,2.841196,0.3550898005272863,0.3550898005272863,3,1,0.5292956666666667,1.587887,0.135651968140022,0.1356519681400219,0.5012775769835556,4.511498192852001,0.04816863030797418,0.04816863030797418),
    @("sCs","Bmp3","Acvr2b","FAPs",3,1,0.9470653333333333,2.841196,0.3550898005272863,0.3550898005272863,3,1,1.362890666666667,4.088672000000001,0.3492921120199358,0.3492921120199358,1.290746503523556,11.616718531712,0.1240300663829135,0.1240300663829136),
    @("sCs","Bmp3","Acvr2b","sCs",3,1,0.9470653333333333,2.841196,0.3550898005272863,0.3550898005272863,3,1,2.009678666666666,6.029036,0.5150559198400423,0.5150559198400423,1.903296996339555,17.129672967056,0.1828911038363986,0.1828911038363986)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
